$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 2, shifting existing data (rows 2-6) down to rows 4-8.
$ws.Rows.Item(2).Resize(2).Insert()

# The inserted rows pick up formatting from the row above (the header row).
# Clear that inherited formatting and reapply the normal data-row format:
# column A keeps the existing date number format, column B stays unstyled.
$ws.Range("A2:B3").ClearFormats()
$ws.Range("A2:A3").NumberFormat = $ws.Range("A4").NumberFormat

# Set the date values (column A) for all data rows 2-8.
$ws.Range("A2").Value = 45046
$ws.Range("A3").Value = 45077
$ws.Range("A4").Value = 45107
$ws.Range("A5").Value = 45138
$ws.Range("A6").Value = 45169
$ws.Range("A7").Value = 45199
$ws.Range("A8").Value = 45230

# Set the RH values (column B) for all data rows 2-8.
$ws.Range("B2").Value = 549500
$ws.Range("B3").Value = 1399000
$ws.Range("B4").Value = 1736000
$ws.Range("B5").Value = 1809500
$ws.Range("B6").Value = 2026000
$ws.Range("B7").Value = 1887500
$ws.Range("B8").Value = 1953000
